$d = $word.ActiveDocument

# wdReplaceOne = 1, wdFindStop (no wrap) = 0

# --- Body text (document.xml): the bold "TERE" salutation run -> "QWER" ---
$d.Content.Find.Execute("TERE", $true, $false, $false, $false, $false, `
    $true, 0, $false, "QWER", 1)

# --- Header text (header1.xml) ---
$hdr = $d.Sections.Item(1).Headers.Item(1)

# "DIRETORIA DE ENSINO REGIAO TRE" -> "...QWER"
$hdr.Range.Find.Execute("TRE", $true, $false, $false, $false, $false, `
    $true, 0, $false, "QWER", 1)

# "TERE - DEP." -> "QWER - DEP."
$hdr.Range.Find.Execute("TERE", $true, $false, $false, $false, $false, `
    $true, 0, $false, "QWER", 1)

# "Tre, n. Tre - Tre - Tre - Tre" -> "Qwer, n. Qwer - Qewr - Qewr - Qwer"
$hdr.Range.Find.Execute("Tre", $true, $false, $false, $false, $false, `
    $true, 0, $false, "Qwer", 1)
$hdr.Range.Find.Execute("Tre", $true, $false, $false, $false, $false, `
    $true, 0, $false, "Qwer", 1)
$hdr.Range.Find.Execute("Tre", $true, $false, $false, $false, $false, `
    $true, 0, $false, "Qewr", 1)
$hdr.Range.Find.Execute("Tre", $true, $false, $false, $false, $false, `
    $true, 0, $false, "Qewr", 1)
$hdr.Range.Find.Execute("Tre", $true, $false, $false, $false, $false, `
    $true, 0, $false, "Qwer", 1)

# "CEP: tre" / "Tel: tre" / "Email: tre" -> "qwer"
$hdr.Range.Find.Execute("tre", $true, $false, $false, $false, $false, `
    $true, 0, $false, "qwer", 1)
$hdr.Range.Find.Execute("tre", $true, $false, $false, $false, $false, `
    $true, 0, $false, "qwer", 1)
$hdr.Range.Find.Execute("tre", $true, $false, $false, $false, $false, `
    $true, 0, $false, "qwer", 1)
